# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with the default (unstyled) format, used to strip the
# "number stored as text" quote-prefix style that Excel applies when a
# numeric-looking string is typed into a cell via a leading apostrophe.
$formatDonor = $ws.Range("D4")

$ws.Range("D2").Value = "'30.316.49"
$formatDonor.Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "'1.871.12"
$formatDonor.Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'243.46"
$formatDonor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  -1.73%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.4695"
$formatDonor.Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = "  -1.02%  "

$ws.Range("E8").Value = "  -0.96%  "

$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").Value = "'22.07"
$formatDonor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("D11").Value = "'0.07765"
$formatDonor.Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "'1.875.52"
$formatDonor.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "'95.56"
$formatDonor.Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").Value = "'0.7198"
$formatDonor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("D16").Value = "'279.31"
$formatDonor.Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("D17").Value = "'30.306.68"
$formatDonor.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").Value = "'0.000007427"
$formatDonor.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").Value = "'2.120.40"
$formatDonor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "'1.000"
$formatDonor.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "'5.223"
$formatDonor.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").Value = "'6.224"
$formatDonor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").Value = "'163.28"
$formatDonor.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").Value = "'9.049"
$formatDonor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").Value = "'18.65"
$formatDonor.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").Value = "'1.878"
$formatDonor.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("D29").Value = "'1.318"
$formatDonor.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "  -2.13%  "

$ws.Range("D30").Value = "'0.09562"
$formatDonor.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  -3.28%  "

$ws.Range("D31").Value = "'1.468"
$formatDonor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  -2.66%  "

$ws.Range("D32").Value = "'4.214"
$formatDonor.Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("D33").Value = "'4.090"
$formatDonor.Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").Value = "'0.04806"
$formatDonor.Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "'0.6872"
$formatDonor.Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").Value = "'2.706"
$formatDonor.Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("D39").Value = "'2.812"
$formatDonor.Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("D40").Value = "'6.220"
$formatDonor.Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "'74.15"
$formatDonor.Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("D42").Value = "'0.4233"
$formatDonor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  +1.72%  "

$ws.Range("D43").Value = "'1.933"
$formatDonor.Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("D44").Value = "'0.9994"
$formatDonor.Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "'0.8249"
$formatDonor.Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  -1.13%  "

$ws.Range("D46").Value = "'100.77"
$formatDonor.Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").Value = "'9.553"
$formatDonor.Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = "  +2.29%  "

$ws.Range("E48").Value = "  -0.67%  "

$ws.Range("D49").Value = "'6.910"
$formatDonor.Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("D50").Value = "'898.73"
$formatDonor.Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").Value = "'0.05718"
$formatDonor.Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  +0.85%  "

$excel.CutCopyMode = 0
